$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Camote, Vega Modelo de Temuco) was added.
# It belongs chronologically among the existing rows, so insert a fresh
# row at 46 (pushing the former rows 46-59 down to 47-60) and populate it.
$ws.Rows.Item(46).Insert()

$ws.Range("A46").Value = 10
$ws.Range("B46").Value = "Vega Modelo de Temuco"
$ws.Range("C46").Value = "La Araucanía"
$ws.Range("D46").Value = 44636
$ws.Range("E46").Value = 9
$ws.Range("F46").Value = 100114002
$ws.Range("G46").Value = "Camote"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 50
$ws.Range("K46").Value = 18000
$ws.Range("L46").Value = 18000
$ws.Range("M46").Value = 18000
$ws.Range("N46").Value = "`$/malla 20 kilos"
$ws.Range("O46").Value = "Perú"
$ws.Range("P46").Value = 900
$ws.Range("Q46").Value = 20
$ws.Range("R46").Value = "Hortaliza"
